# Applies the "mapping.xlsx" revision described by the commit diff:
#  - Expands Table1 (Account/Account Name/Task) from A1:C61 to A1:C67,
#    inserting several new account rows (and renumbering 7700/7701/.. order).
#  - Expands Table2 (Projects/Spec4) from E1:F17 to E1:F19, adding two rows.
#  - Adds the new strings this introduces to the shared-strings table.
#  - Updates the active cell selection left behind by the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Table1 (A:C) - rows 59..61 get replaced/reordered and rows 62..67 are
# brand-new, as account codes 7700-7730 (plus a late-added 7320 row) are
# filled in around the existing 7701/7704/7710 entries.
# ---------------------------------------------------------------------
$table1Rows = @(
    @{ Row = 59; Account = 7700; Name = "Interne avdelingsseminar";       Task = 99 },
    @{ Row = 60; Account = 7701; Name = "Møtekostander";                  Task = 99 },
    @{ Row = 61; Account = 7702; Name = "Møtekostnader utlandet";         Task = 99 },
    @{ Row = 62; Account = 7703; Name = "Div interne prosjektkostnader";  Task = 99 },
    @{ Row = 63; Account = 7704; Name = "Konferanser deltakelse";         Task = 99 },
    @{ Row = 64; Account = 7710; Name = "Styrekostnader";                 Task = 80 },
    @{ Row = 65; Account = 7713; Name = "Prosjektkostnader Towards2040";  Task = 99 },
    @{ Row = 66; Account = 7730; Name = "Medieovervåkning / presseklipp"; Task = 80 },
    @{ Row = 67; Account = 7320; Name = "Markering - arrangementer";      Task = 80 }
)

# Grow the table first so the new rows become part of Table1 (and its
# AutoFilter) rather than plain out-of-table values.
$table1 = $ws.ListObjects.Item("Table1")
$table1.Resize($ws.Range("A1:C67"))

foreach ($r in $table1Rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Account
    $ws.Cells.Item($r.Row, 2).Value = $r.Name
    $ws.Cells.Item($r.Row, 3).Value = $r.Task
}

# ---------------------------------------------------------------------
# Table2 (E:F) - two new rows appended (E18:F18 and E19:F19).
# E column keeps its existing Text number format; E18 holds a genuine
# number (31316) while E19 holds a text value ("31279").
# ---------------------------------------------------------------------
$table2 = $ws.ListObjects.Item("Table2")
$table2.Resize($ws.Range("E1:F19"))

$ws.Range("E18").Value = 31316
$ws.Range("E18").NumberFormat = "@"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "Towards2040"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "31279"
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "Oppdrag"

# ---------------------------------------------------------------------
# Leftover UI state from the edit: the cursor ended up on I16.
# ---------------------------------------------------------------------
$ws.Range("I16").Select()
